$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Variável" (B) column text for all data rows (2-10): 2024/07 -> 2024/10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Diferença 2024/10 - 2024/10"
}

# Update region names (column A) where they changed
$ws.Range("A3").Value  = "Rio Grande do Norte"
$ws.Range("A4").Value  = "Maranhão"
$ws.Range("A5").Value  = "Distrito Federal"
$ws.Range("A6").Value  = "Amazonas"
$ws.Range("A7").Value  = "Roraima"

# Update values (column C)
$ws.Range("C2").Value  = 1.09
$ws.Range("C3").Value  = 1.02
$ws.Range("C4").Value  = 0.97
$ws.Range("C5").Value  = 0.95
$ws.Range("C7").Value  = 0.9399999999999999
$ws.Range("C8").Value  = 0.75
$ws.Range("C9").Value  = 0.84
$ws.Range("C10").Value = 0.83

# Update placement (column D) for row 8
$ws.Range("D8").Value = "20º"
